$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared-string / header text values
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Update column C (GDP) values for rows 2-27
$ws.Range("C2").Value = 7772.38875590225
$ws.Range("C3").Value = 9502.243585046588
$ws.Range("C4").Value = 1909.084588129339
$ws.Range("C5").Value = 6128.19547247793
$ws.Range("C6").Value = 4547.50930098406
$ws.Range("C7").Value = 4729.735976516416
$ws.Range("C8").Value = 8082.02845866252
$ws.Range("C9").Value = 10385.96443195552
$ws.Range("C10").Value = 1955.461557360978
$ws.Range("C11").Value = 6336.709213679884
$ws.Range("C12").Value = 4633.590358399045
$ws.Range("C13").Value = 5082.354756663512
$ws.Range("C14").Value = 8841.561277324312
$ws.Range("C15").Value = 2024.117324382548
$ws.Range("C16").Value = 6711.616186806423
$ws.Range("C17").Value = 4921.848409120176
$ws.Range("C18").Value = 5360.226632400601
$ws.Range("C19").Value = 2094.024217383061
$ws.Range("C20").Value = 2201.396847776877
$ws.Range("C21").Value = 5996.49696468919
$ws.Range("C22").Value = 6114.227214287786
$ws.Range("C23").Value = 3382.563653843273
$ws.Range("C24").Value = 514.0573067519859
$ws.Range("C25").Value = 1875.732161108182
$ws.Range("C26").Value = 1895.214690888655
$ws.Range("C27").Value = 0
